$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 512; this shifts rows 512..619 down to 513..620,
# carrying along their existing formatting (mirrors the canonical diff, which is
# equivalent to one new pricing record being inserted before the former row 512).
$ws.Rows.Item(512).EntireRow.Insert()

# Populate the newly inserted row 512 with the new record's data.
$ws.Range("A512").Value = 7
$ws.Range("B512").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C512").Value = "Ñuble"
$ws.Range("D512").Value = 45244
$ws.Range("E512").Value = 16
$ws.Range("F512").Value = "Fruta"
$ws.Range("G512").Value = 100104
$ws.Range("H512").Value = "Frutos de pepita"
$ws.Range("I512").Value = 100104005
$ws.Range("J512").Value = "Pera"
$ws.Range("K512").Value = "Packham's Triumph"
$ws.Range("L512").Value = "Primera"
$ws.Range("M512").Value = 180
$ws.Range("N512").Value = 15000
$ws.Range("O512").Value = 16000
$ws.Range("P512").Value = 15444
$ws.Range("Q512").Value = "$/bandeja 18 kilos granel"
$ws.Range("R512").Value = "Región de O'Higgins"
$ws.Range("S512").Value = 858
$ws.Range("T512").Value = 18
